# Generate Report for Handback
#
# Removes the "9fed5074-b05e-4132-85b7-1f44abe48eba" file's row (row 3)
# from every sheet (Overview, zh-cn, de-de), shifting the
# ".localization-config" row up from row 4 to row 3, fixes up the
# hyperlinks that referenced the now-deleted / now-shifted cells, and
# bumps the handoff/handback timestamps recorded for the
# "36491f5a-d66c-495f-9f55-eaba4cdc0280" entry (new handback run).

function Remove-HyperlinkAt($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            break
        }
    }
}

$wb = $excel.ActiveWorkbook

$localizationConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/320bb16eaa9fe3e04694c4debc6f1913119bb47d/.localization-config"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop the hyperlink that lives on the row we are about to remove.
Remove-HyperlinkAt $wsOverview '$A$3'
# Drop the hyperlink for ".localization-config" too - it will be
# re-added at its new location once the row shift has happened.
Remove-HyperlinkAt $wsOverview '$A$4'

# Remove the whole "9fed5074-...md" row - everything below shifts up.
$wsOverview.Rows(3).Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $localizationConfigUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Remove-HyperlinkAt $wsZhCn '$A$3'
Remove-HyperlinkAt $wsZhCn '$C$3'
Remove-HyperlinkAt $wsZhCn '$E$3'
Remove-HyperlinkAt $wsZhCn '$F$3'
Remove-HyperlinkAt $wsZhCn '$A$4'

$wsZhCn.Rows(3).Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $localizationConfigUrl, "", "", ".localization-config") | Out-Null

# New handback run for the 36491f5a... zh-cn handoff/handback timestamps.
$wsZhCn.Range("D2").Value = "2016-03-09 03:15:56"
$wsZhCn.Range("G2").Value = "2016-03-09 03:16:53"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Remove-HyperlinkAt $wsDeDe '$A$3'
Remove-HyperlinkAt $wsDeDe '$C$3'
Remove-HyperlinkAt $wsDeDe '$E$3'
Remove-HyperlinkAt $wsDeDe '$F$3'
Remove-HyperlinkAt $wsDeDe '$A$4'

$wsDeDe.Rows(3).Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $localizationConfigUrl, "", "", ".localization-config") | Out-Null

# New handback run for the 36491f5a... de-de handoff/handback timestamps.
$wsDeDe.Range("D2").Value = "2016-03-09 03:16:11"
$wsDeDe.Range("G2").Value = "2016-03-09 03:17:29"
